$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H58").Value = 333.07144
$ws.Range("I58").Value = 204.84616
$ws.Range("J58").Value = 2000
$ws.Range("K58").Value = 614.5384799999999
$ws.Range("L58").Value = 6000
$ws.Range("M58").Value = -464.5384799999999
$ws.Range("N58").Value = -6300

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H94").Value = 1500
$ws.Range("I94").Value = 1500
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 1500
$ws.Range("L94").Value = 0
$ws.Range("M94").Value = -1049
$ws.Range("N94").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 1401
$ws.Range("I100").Value = 1351.25
$ws.Range("J100").Value = 1600
$ws.Range("K100").Value = 1351.25
$ws.Range("L100").Value = 1600
$ws.Range("M100").Value = -810.25
$ws.Range("N100").Value = -2682

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 1498.8077
$ws.Range("J112").Value = 1659.3889
$ws.Range("L112").Value = 4978.1667
$ws.Range("N112").Value = -7194.1667

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 10160.954
$ws.Range("I116").Value = 2589.0833
$ws.Range("J116").Value = 19247.2
$ws.Range("K116").Value = 2589.0833
$ws.Range("L116").Value = 19247.2
$ws.Range("M116").Value = 852.9167000000002
$ws.Range("N116").Value = -26131.2

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H130").Value = 33631
$ws.Range("J130").Value = 33631
$ws.Range("L130").Value = 33631
$ws.Range("N130").Value = -43671

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 4247.6206
$ws.Range("I94").Value = 517.4375
$ws.Range("J94").Value = 8838.615
$ws.Range("K94").Value = 517.4375
$ws.Range("L94").Value = 8838.615
$ws.Range("M94").Value = -66.4375
$ws.Range("N94").Value = -9740.615

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H96").Value = 20314.834
$ws.Range("I96").Value = 7614
$ws.Range("J96").Value = 26665.25
$ws.Range("K96").Value = 7614
$ws.Range("L96").Value = 26665.25
$ws.Range("M96").Value = -4868
$ws.Range("N96").Value = -32157.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1463.1666
$ws.Range("I99").Value = 1066.3334
$ws.Range("J99").Value = 1860
$ws.Range("K99").Value = 1066.3334
$ws.Range("L99").Value = 1860
$ws.Range("M99").Value = 431.6666
$ws.Range("N99").Value = -4856

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 34.61111
$ws.Range("I7").Value = 24.428572
$ws.Range("J7").Value = 70.25
$ws.Range("K7").Value = 24.428572
$ws.Range("L7").Value = 70.25
$ws.Range("M7").Value = 88.571428
$ws.Range("N7").Value = -296.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3888.3447
$ws.Range("I31").Value = 1566
$ws.Range("J31").Value = 4333.726
$ws.Range("K31").Value = 1566
$ws.Range("L31").Value = 4333.726
$ws.Range("M31").Value = -1271
$ws.Range("N31").Value = -4923.726

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 3888.3447
$ws.Range("I34").Value = 1566
$ws.Range("J34").Value = 4333.726
$ws.Range("K34").Value = 1566
$ws.Range("L34").Value = 4333.726
$ws.Range("M34").Value = -1364
$ws.Range("N34").Value = -4737.726

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 1406
$ws.Range("I99").Value = 1271.8096
$ws.Range("J99").Value = 1687.8
$ws.Range("K99").Value = 1271.8096
$ws.Range("L99").Value = 1687.8
$ws.Range("M99").Value = 226.1904
$ws.Range("N99").Value = -4683.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H100").Value = 22319.75
$ws.Range("J100").Value = 22319.75
$ws.Range("L100").Value = 22319.75
$ws.Range("N100").Value = -24483.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H106").Value = 20079.5
$ws.Range("J106").Value = 20079.5
$ws.Range("L106").Value = 20079.5
$ws.Range("N106").Value = -22603.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 1406
$ws.Range("I126").Value = 1271.8096
$ws.Range("J126").Value = 1687.8
$ws.Range("K126").Value = 3815.4288
$ws.Range("L126").Value = 5063.4
$ws.Range("M126").Value = -1345.4288
$ws.Range("N126").Value = -10003.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H13").Value = 998
$ws.Range("I13").Value = 230
$ws.Range("J13").Value = 3302
$ws.Range("K13").Value = 690
$ws.Range("L13").Value = 9906
$ws.Range("M13").Value = -522
$ws.Range("N13").Value = -10242

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H109").Value = 1765.2727
$ws.Range("I109").Value = 445.42856
$ws.Range("J109").Value = 4075
$ws.Range("K109").Value = 1336.28568
$ws.Range("L109").Value = 12225
$ws.Range("M109").Value = -296.28568
$ws.Range("N109").Value = -14305

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 606422.06
$ws.Range("I113").Value = 632.4211
$ws.Range("J113").Value = 1181922.2
$ws.Range("K113").Value = 1897.2633
$ws.Range("L113").Value = 3545766.6
$ws.Range("M113").Value = 272.7366999999999
$ws.Range("N113").Value = -3550106.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 1365.7142
$ws.Range("I129").Value = 1765
$ws.Range("J129").Value = 1299.1666
$ws.Range("K129").Value = 5295
$ws.Range("L129").Value = 3897.4998
$ws.Range("M129").Value = -295
$ws.Range("N129").Value = -13897.4998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 784.81177
$ws.Range("J131").Value = 933.9394
$ws.Range("L131").Value = 2801.8182
$ws.Range("N131").Value = -12881.8182

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H24").Value = 95307.695
$ws.Range("I24").Value = 121700
$ws.Range("J24").Value = 7333.3335
$ws.Range("K24").Value = 121700
$ws.Range("L24").Value = 7333.3335
$ws.Range("M24").Value = -121527
$ws.Range("N24").Value = -7679.3335

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H92").Value = 32000
$ws.Range("J92").Value = 32000
$ws.Range("L92").Value = 32000
$ws.Range("N92").Value = -36992

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 2354.2856
$ws.Range("I100").Value = 2357.5
$ws.Range("J100").Value = 2350
$ws.Range("K100").Value = 2357.5
$ws.Range("L100").Value = 2350
$ws.Range("M100").Value = -1816.5
$ws.Range("N100").Value = -3432

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H104").Value = 16040
$ws.Range("J104").Value = 16040
$ws.Range("L104").Value = 16040
$ws.Range("N104").Value = -23028

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 62450.766
$ws.Range("I122").Value = 114011.445
$ws.Range("K122").Value = 342034.335
$ws.Range("M122").Value = -339584.335

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1743.3334
$ws.Range("I96").Value = 1000
$ws.Range("J96").Value = 1892
$ws.Range("K96").Value = 1000
$ws.Range("L96").Value = 1892
$ws.Range("M96").Value = 373
$ws.Range("N96").Value = -4638

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H99").Value = 35738
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 35738
$ws.Range("K99").Value = 0
$ws.Range("L99").Value = 35738
$ws.Range("M99").ClearContents()
$ws.Range("N99").Value = -41728

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H104").Value = 16175
$ws.Range("J104").Value = 16175
$ws.Range("L104").Value = 16175
$ws.Range("N104").Value = -23163
